$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.133054
$ws.Range("H2").Value = 0.399162
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4660823333333333
$ws.Range("N2").Value = 1.398247
$ws.Range("O2").Value = 0.02470419871925173
$ws.Range("P2").Value = 0.02470419871925173
$ws.Range("Q2").Value = 0.06201411877933333
$ws.Range("R2").Value = 0.558127069014
$ws.Range("S2").Value = 0.02470419871925173
$ws.Range("T2").Value = 0.02470419871925173

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.133054
$ws.Range("H3").Value = 0.399162
$ws.Range("O3").Value = 0.6873021241557511
$ws.Range("P3").Value = 0.6873021241557511
$ws.Range("Q3").Value = 1.725311395405333
$ws.Range("R3").Value = 15.527802558648
$ws.Range("S3").Value = 0.6873021241557511
$ws.Range("T3").Value = 0.6873021241557511

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.133054
$ws.Range("H4").Value = 0.399162
$ws.Range("O4").Value = 0.2879936771249972
$ws.Range("P4").Value = 0.2879936771249972
$ws.Range("Q4").Value = 0.7229408370573334
$ws.Range("R4").Value = 6.506467533516001
$ws.Range("S4").Value = 0.2879936771249972
$ws.Range("T4").Value = 0.2879936771249972
